# Update generated output numbers (想去人数 / 最低票价) across sheets
# 展览 (sheet1), 演出 (sheet2), 全部类型 (sheet4)

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 4952
$ws1.Range("F4").Value  = 66
$ws1.Range("F5").Value  = 2845
$ws1.Range("F9").Value  = 1748
$ws1.Range("G9").Value  = 70
$ws1.Range("F17").Value = 97
$ws1.Range("F19").Value = 1053
$ws1.Range("F22").Value = 686
$ws1.Range("F23").Value = 762
$ws1.Range("F24").Value = 158
$ws1.Range("F26").Value = 45
$ws1.Range("F27").Value = 573
$ws1.Range("F28").Value = 65
$ws1.Range("F30").Value = 1713
$ws1.Range("F31").Value = 421
$ws1.Range("F33").Value = 1614
$ws1.Range("F34").Value = 231
$ws1.Range("F35").Value = 2429
$ws1.Range("F36").Value = 425
$ws1.Range("F39").Value = 122
$ws1.Range("F42").Value = 831
$ws1.Range("F43").Value = 1532
$ws1.Range("F44").Value = 241
$ws1.Range("F47").Value = 80
$ws1.Range("F48").Value = 85
$ws1.Range("F49").Value = 123

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F12").Value = 50

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 4952
$ws4.Range("F4").Value  = 2845
$ws4.Range("F5").Value  = 1748
$ws4.Range("G5").Value  = 70
$ws4.Range("F14").Value = 97
$ws4.Range("F15").Value = 1053
$ws4.Range("F17").Value = 686
$ws4.Range("F18").Value = 762
$ws4.Range("F19").Value = 158
$ws4.Range("F25").Value = 45
$ws4.Range("F26").Value = 573
$ws4.Range("F28").Value = 1713
$ws4.Range("F29").Value = 422
$ws4.Range("F33").Value = 2429
$ws4.Range("F34").Value = 425
$ws4.Range("F39").Value = 50
$ws4.Range("F40").Value = 122
$ws4.Range("F43").Value = 831
$ws4.Range("F44").Value = 1532
$ws4.Range("F46").Value = 241
$ws4.Range("F48").Value = 80
$ws4.Range("F49").Value = 85
